$d = $word.ActiveDocument

# Replace the two occurrences of "Richy Rich" with "%name%"
$d.Content.Find.Execute("Richy Rich", $true, $false, $false, $false, $false,
                         $true, 1, $false, "%name%", 2)

# Replace job title placeholder
$d.Content.Find.Execute("this is the job title", $true, $false, $false, $false, $false,
                         $true, 1, $false, "%jobtitle%", 2)

# Replace phone number
$d.Content.Find.Execute("513-867-5309", $true, $false, $false, $false, $false,
                         $true, 1, $false, "%phone%", 2)

# Replace email
$d.Content.Find.Execute("email@mail.uc.edu", $true, $false, $false, $false, $false,
                         $true, 1, $false, "%email%", 2)

# Replace website
$d.Content.Find.Execute("https://github.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "%website%", 2)

# Replace "Dear Hiring Manager" with "Dear %recipient%"
$d.Content.Find.Execute("Dear Hiring Manager", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dear %recipient%", 2)

# Replace body placeholder
$d.Content.Find.Execute("this is the body", $true, $false, $false, $false, $false,
                         $true, 1, $false, "%body%", 2)
